$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "26.982.90"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").Value = "1.881.75"
$ws.Range("E3").Value = "  +1.39%  "

# Row 4
Set-TextValue "D4" "0.9973"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
Set-TextValue "D5" "306.65"
$ws.Range("E5").Value = "  +0.67%  "

# Row 6
Set-TextValue "D6" "0.9982"
$ws.Range("E6").Value = "  -0.21%  "

# Row 7
Set-TextValue "D7" "0.5173"
$ws.Range("E7").Value = "  +1.62%  "

# Row 8
Set-TextValue "D8" "0.3716"
$ws.Range("E8").Value = "  +1.83%  "

# Row 9
Set-TextValue "D9" "0.07192"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
Set-TextValue "D10" "0.9013"
$ws.Range("E10").Value = "  +1.43%  "

# Row 11
Set-TextValue "D11" "20.81"
$ws.Range("E11").Value = "  +0.34%  "

# Row 12
Set-TextValue "D12" "0.07564"
$ws.Range("E12").Value = "  +0.48%  "

# Row 13
$ws.Range("D13").Value = "1.892.63"
$ws.Range("E13").Value = "  +2.00%  "

# Row 14
Set-TextValue "D14" "95.12"
$ws.Range("E14").Value = "  +3.90%  "

# Row 15
Set-TextValue "D15" "5.254"
$ws.Range("E15").Value = "  +0.37%  "

# Row 16
Set-TextValue "D16" "0.9990"
$ws.Range("E16").Value = "  -0.13%  "

# Row 17
Set-TextValue "D17" "0.000008489"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
Set-TextValue "D18" "14.25"
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
Set-TextValue "D19" "0.9976"
$ws.Range("E19").Value = "  -0.30%  "

# Row 20
$ws.Range("D20").Value = "26.999.08"
$ws.Range("E20").Value = "  +0.54%  "

# Row 21
Set-TextValue "D21" "5.036"
$ws.Range("E21").Value = "  +0.48%  "

# Row 22
$ws.Range("D22").Value = "2.112.99"
$ws.Range("E22").Value = "  +1.30%  "

# Row 23
Set-TextValue "D23" "10.41"
$ws.Range("E23").Value = "  +1.49%  "

# Row 24
Set-TextValue "D24" "6.451"
$ws.Range("E24").Value = "  +0.11%  "

# Row 25
Set-TextValue "D25" "145.71"
$ws.Range("E25").Value = "  -0.53%  "

# Row 26
Set-TextValue "D26" "1.783"
$ws.Range("E26").Value = "  -2.11%  "

# Row 27
$ws.Range("E27").Value = "  +0.80%  "

# Row 28
Set-TextValue "D28" "2.118"
$ws.Range("E28").Value = "  +3.04%  "

# Row 29
Set-TextValue "D29" "114.49"
$ws.Range("E29").Value = "  +1.26%  "

# Row 30
Set-TextValue "D30" "4.962"
$ws.Range("E30").Value = "  +6.09%  "

# Row 31
Set-TextValue "D31" "4.772"
$ws.Range("E31").Value = "  +2.90%  "

# Row 32
Set-TextValue "D32" "0.09200"
$ws.Range("E32").Value = "  -0.71%  "

# Row 33
Set-TextValue "D33" "0.05033"
$ws.Range("E33").Value = "  -1.62%  "

# Row 34
Set-TextValue "D34" "0.7581"
$ws.Range("E34").Value = "  +3.64%  "

# Row 35
Set-TextValue "D35" "1.181"
$ws.Range("E35").Value = "  +2.61%  "

# Row 36
Set-TextValue "D36" "3.003"
$ws.Range("E36").Value = "  -2.11%  "

# Row 37
Set-TextValue "D37" "3.280"
$ws.Range("E37").Value = "  +2.93%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "2.498"
$ws.Range("E38").Value = "  +1.67%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.01992"
$ws.Range("E39").Value = "  -0.71%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D40" "0.5585"
$ws.Range("E40").Value = "  +5.70%  "

# Row 41
Set-TextValue "D41" "1.075"
$ws.Range("E41").Value = "  +0.24%  "

# Row 42
Set-TextValue "D42" "6.585"
$ws.Range("E42").Value = "  +1.52%  "

# Row 43
Set-TextValue "D43" "117.17"
$ws.Range("E43").Value = "  -0.24%  "

# Row 44
Set-TextValue "D44" "8.808"
$ws.Range("E44").Value = "  +4.81%  "

# Row 45
Set-TextValue "D45" "0.1503"
$ws.Range("E45").Value = "  +2.16%  "

# Row 46
Set-TextValue "D46" "0.4784"
$ws.Range("E46").Value = "  +3.10%  "

# Row 47
Set-TextValue "D47" "10.17"
$ws.Range("E47").Value = "  +2.23%  "

# Row 48
Set-TextValue "D48" "0.9982"
$ws.Range("E48").Value = "  -0.20%  "

# Row 49
Set-TextValue "D49" "1.568"
$ws.Range("E49").Value = "  +0.58%  "

# Row 50
Set-TextValue "D50" "37.08"
$ws.Range("E50").Value = "  +0.19%  "

# Row 51
Set-TextValue "D51" "63.40"
$ws.Range("E51").Value = "  +0.68%  "
